$d = $word.ActiveDocument

# --- 1. Collapse the split runs that make up {{RECIP_ADR2}} and {{RECIP_ADR3}} into a
#        single run each (no textual change, just simplifies the run structure) ---
$adr2Para = $d.Paragraphs.Item(11)
if ($adr2Para.Range.Text.Trim() -ne "{{RECIP_ADR2}}") { throw "Unexpected paragraph 11: $($adr2Para.Range.Text)" }
$adr2Para.Range.Find.Execute("{{RECIP_ADR2}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{RECIP_ADR2}}", 2) | Out-Null

$adr3Para = $d.Paragraphs.Item(12)
if ($adr3Para.Range.Text.Trim() -ne "{{RECIP_ADR3}}") { throw "Unexpected paragraph 12: $($adr3Para.Range.Text)" }
$adr3Para.Range.Find.Execute("{{RECIP_ADR3}}", $true, $false, $false, $false, $false, $true, 1, $false, "{{RECIP_ADR3}}", 2) | Out-Null

# --- 2. After the {{RECIP_ADR3}} paragraph, add a blank paragraph followed by a
#        centered "Job Announcement: {{ANNOUNCEMENT_NUM}}" paragraph ---
$d.Paragraphs.Item(12).Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs.Item(13).Range.InsertParagraphAfter() | Out-Null
$jobPara = $d.Paragraphs.Item(14)
$jobPara.Range.Text = "Job Announcement: {{ANNOUNCEMENT_NUM}}"
$jobPara.Alignment = 1

# --- 3. After "Dear {RECIP_NAME}," add a blank paragraph, a body paragraph that
#        introduces the applicant, and a paragraph holding three spaces ---
$dearPara = $d.Paragraphs.Item(16)
if ($dearPara.Range.Text.Trim() -ne "Dear {RECIP_NAME},") { throw "Unexpected paragraph 16: $($dearPara.Range.Text)" }

$dearPara.Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs.Item(17).Range.InsertParagraphAfter() | Out-Null
$bodyPara = $d.Paragraphs.Item(18)
$bodyRange = $bodyPara.Range
$bodyRange.End = $bodyRange.End - 1
$bodyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">I am contacting you today to express my interest in the {{JOB_POS}} position.  My specialized experience is closely related to the job description, and I believe that my knowledge, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>skills</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and abilities can greatly aid your organization.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$bodyRange.InsertXML($bodyXml) | Out-Null
$d.Paragraphs.Item(18).SpaceAfter = 0

$d.Paragraphs.Item(18).Range.InsertParagraphAfter() | Out-Null
$spacesPara = $d.Paragraphs.Item(19)
$spacesPara.Range.Text = "   "
